# Update "想去人数" (number of people interested) figures on the
# 展览 and 全部类型 sheets, which share the same underlying data rows.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F20" = 3196
    "F21" = 70
    "F32" = 572
    "F33" = 1823
    "F34" = 278
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
